# Localize the built-in slide-layout display names (and a couple of
# placeholder names/prompt text) from English to Turkish.
#
# Each of the 11 slide layouts is reached through
# $p.SlideMaster.CustomLayouts.Item(n) (n is 1-based, matching
# slideLayout<n>.xml) and its cSld/@name is exposed as the normal
# PowerPoint COM `.Name` property on the CustomLayout object.

$p = $ppt.ActivePresentation
$layouts = $p.SlideMaster.CustomLayouts

# slideLayout1.xml : "Title Slide" -> "Başlık Slaydı"
$layouts.Item(1).Name = "Başlık Slaydı"

# slideLayout2.xml : "Title and Content" -> "Başlık ve İçerik"
$layouts.Item(2).Name = "Başlık ve İçerik"

# slideLayout3.xml : "Section Header" -> "Bölüm Üstbilgisi"
$layouts.Item(3).Name = "Bölüm Üstbilgisi"

# slideLayout4.xml : "Two Content" -> "İki İçerik"
#   plus its Title placeholder shape "Title 1" -> "Başlık 1"
$layout4 = $layouts.Item(4)
$layout4.Name = "İki İçerik"
try {
    $layout4.Shapes.Item(1).Name = "Başlık 1"
} catch {
    Write-Output "skip: layout4 shape rename unavailable ($_)"
}

# slideLayout5.xml : "Comparison" -> "Karşılaştırma"
$layouts.Item(5).Name = "Karşılaştırma"

# slideLayout6.xml : "Title Only" -> "Yalnızca Başlık"
$layouts.Item(6).Name = "Yalnızca Başlık"

# slideLayout7.xml : "Blank" -> "Boş"
$layouts.Item(7).Name = "Boş"

# slideLayout8.xml : "Content with Caption" -> "Başlıklı İçerik"
$layouts.Item(8).Name = "Başlıklı İçerik"

# slideLayout9.xml : "Picture with Caption" -> "Başlıklı Resim"
#   plus the picture placeholder prompt text and the date placeholder name
$layout9 = $layouts.Item(9)
$layout9.Name = "Başlıklı Resim"
$layout9.Shapes.Item(2).TextFrame.TextRange.Text = "Resim eklemek için simgeye tıklayın"
try {
    $layout9.Shapes.Item(4).Name = "Veri Yer Tutucusu 4"
} catch {
    Write-Output "skip: layout9 shape rename unavailable ($_)"
}

# slideLayout10.xml : "Title and Vertical Text" -> "Başlık, Dikey Metin"
$layouts.Item(10).Name = "Başlık, Dikey Metin"

# slideLayout11.xml : "Vertical Title and Text" -> "Dikey Başlık ve Metin"
$layouts.Item(11).Name = "Dikey Başlık ve Metin"
